$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.099562666666667
$ws.Range("H2").Value = 9.298688
$ws.Range("I2").Value = 0.2686390288432488
$ws.Range("J2").Value = 0.2686390288432488
$ws.Range("M2").Value = 5.844648666666667
$ws.Range("N2").Value = 17.533946
$ws.Range("O2").Value = 0.3204643139023235
$ws.Range("P2").Value = 0.3204643139023235
$ws.Range("Q2").Value = 18.11585480698311
$ws.Range("R2").Value = 163.042693262848
$ws.Range("S2").Value = 0.08608922206563822
$ws.Range("T2").Value = 0.08608922206563822
$ws.Range("G3").Value = 3.099562666666667
$ws.Range("H3").Value = 9.298688
$ws.Range("I3").Value = 0.2686390288432488
$ws.Range("J3").Value = 0.2686390288432488
$ws.Range("O3").Value = 0.2396231127748354
$ws.Range("P3").Value = 0.2396231127748355
$ws.Range("Q3").Value = 13.54589990556445
$ws.Range("R3").Value = 121.91309915008
$ws.Range("S3").Value = 0.06437212030422806
$ws.Range("T3").Value = 0.06437212030422808
$ws.Range("G4").Value = 3.099562666666667
$ws.Range("H4").Value = 9.298688
$ws.Range("I4").Value = 0.2686390288432488
$ws.Range("J4").Value = 0.2686390288432488
$ws.Range("M4").Value = 6.452372666666666
$ws.Range("N4").Value = 19.357118
$ws.Range("O4").Value = 0.3537860524377295
$ws.Range("P4").Value = 0.3537860524377295
$ws.Range("Q4").Value = 19.99953342902045
$ws.Range("R4").Value = 179.995800861184
$ws.Range("S4").Value = 0.09504074154515832
$ws.Range("T4").Value = 0.09504074154515833
$ws.Range("G5").Value = 3.099562666666667
$ws.Range("H5").Value = 9.298688
$ws.Range("I5").Value = 0.2686390288432488
$ws.Range("J5").Value = 0.2686390288432488
$ws.Range("M5").Value = 1.570781
$ws.Range("N5").Value = 4.712343
$ws.Range("O5").Value = 0.08612652088511148
$ws.Range("P5").Value = 0.0861265208851115
$ws.Range("Q5").Value = 4.868734145109333
$ws.Range("R5").Value = 43.818607305984
$ws.Range("S5").Value = 0.02313694492822413
$ws.Range("T5").Value = 0.02313694492822413
$ws.Range("G6").Value = 6.189892666666666
$ws.Range("I6").Value = 0.5364778626674904
$ws.Range("J6").Value = 0.5364778626674905
$ws.Range("M6").Value = 5.844648666666667
$ws.Range("N6").Value = 17.533946
$ws.Range("O6").Value = 0.3204643139023235
$ws.Range("P6").Value = 0.3204643139023235
$ws.Range("Q6").Value = 36.17774792104311
$ws.Range("R6").Value = 325.599731289388
$ws.Range("S6").Value = 0.1719220101835223
$ws.Range("T6").Value = 0.1719220101835223
$ws.Range("G7").Value = 6.189892666666666
$ws.Range("I7").Value = 0.5364778626674904
$ws.Range("J7").Value = 0.5364778626674905
$ws.Range("O7").Value = 0.2396231127748354
$ws.Range("P7").Value = 0.2396231127748355
$ws.Range("S7").Value = 0.1285524953871748
$ws.Range("T7").Value = 0.1285524953871748
$ws.Range("G8").Value = 6.189892666666666
$ws.Range("I8").Value = 0.5364778626674904
$ws.Range("J8").Value = 0.5364778626674905
$ws.Range("M8").Value = 6.452372666666666
$ws.Range("N8").Value = 19.357118
$ws.Range("O8").Value = 0.3537860524377295
$ws.Range("P8").Value = 0.3537860524377295
$ws.Range("Q8").Value = 39.93949425200044
$ws.Range("R8").Value = 359.455448268004
$ws.Range("S8").Value = 0.1897983852533618
$ws.Range("T8").Value = 0.1897983852533619
$ws.Range("G9").Value = 6.189892666666666
$ws.Range("I9").Value = 0.5364778626674904
$ws.Range("J9").Value = 0.5364778626674905
$ws.Range("M9").Value = 1.570781
$ws.Range("N9").Value = 4.712343
$ws.Range("O9").Value = 0.08612652088511148
$ws.Range("P9").Value = 0.0861265208851115
$ws.Range("Q9").Value = 9.722965792839332
$ws.Range("R9").Value = 87.50669213555399
$ws.Range("S9").Value = 0.04620497184343159
$ws.Range("T9").Value = 0.0462049718434316
$ws.Range("G10").Value = 1.888584
$ws.Range("H10").Value = 5.665752
$ws.Range("I10").Value = 0.1636835341659699
$ws.Range("J10").Value = 0.1636835341659699
$ws.Range("M10").Value = 5.844648666666667
$ws.Range("N10").Value = 17.533946
$ws.Range("O10").Value = 0.3204643139023235
$ws.Range("P10").Value = 0.3204643139023235
$ws.Range("Q10").Value = 11.038109957488
$ws.Range("R10").Value = 99.34298961739201
$ws.Range("S10").Value = 0.05245473147360508
$ws.Range("T10").Value = 0.05245473147360509
$ws.Range("G11").Value = 1.888584
$ws.Range("H11").Value = 5.665752
$ws.Range("I11").Value = 0.1636835341659699
$ws.Range("J11").Value = 0.1636835341659699
$ws.Range("O11").Value = 0.2396231127748354
$ws.Range("P11").Value = 0.2396231127748355
$ws.Range("Q11").Value = 8.253606259480001
$ws.Range("R11").Value = 74.28245633532001
$ws.Range("S11").Value = 0.03922235796683583
$ws.Range("T11").Value = 0.03922235796683585
$ws.Range("G12").Value = 1.888584
$ws.Range("H12").Value = 5.665752
$ws.Range("I12").Value = 0.1636835341659699
$ws.Range("J12").Value = 0.1636835341659699
$ws.Range("M12").Value = 6.452372666666666
$ws.Range("N12").Value = 19.357118
$ws.Range("O12").Value = 0.3537860524377295
$ws.Range("P12").Value = 0.3537860524377295
$ws.Range("Q12").Value = 12.185847780304
$ws.Range("R12").Value = 109.672630022736
$ws.Range("S12").Value = 0.05790895140163471
$ws.Range("T12").Value = 0.05790895140163473
$ws.Range("G13").Value = 1.888584
$ws.Range("H13").Value = 5.665752
$ws.Range("I13").Value = 0.1636835341659699
$ws.Range("J13").Value = 0.1636835341659699
$ws.Range("M13").Value = 1.570781
$ws.Range("N13").Value = 4.712343
$ws.Range("O13").Value = 0.08612652088511148
$ws.Range("P13").Value = 0.0861265208851115
$ws.Range("Q13").Value = 2.966551864104
$ws.Range("R13").Value = 26.698966776936
$ws.Range("S13").Value = 0.01409749332389426
$ws.Range("T13").Value = 0.01409749332389427
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3599813333333333
$ws.Range("H14").Value = 1.079944
$ws.Range("I14").Value = 0.03119957432329092
$ws.Range("J14").Value = 0.03119957432329093
$ws.Range("M14").Value = 5.844648666666667
$ws.Range("N14").Value = 17.533946
$ws.Range("O14").Value = 0.3204643139023235
$ws.Range("P14").Value = 0.3204643139023235
$ws.Range("Q14").Value = 2.103964419891556
$ws.Range("R14").Value = 18.935679779024
$ws.Range("S14").Value = 0.009998350179557976
$ws.Range("T14").Value = 0.009998350179557977
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3599813333333333
$ws.Range("H15").Value = 1.079944
$ws.Range("I15").Value = 0.03119957432329092
$ws.Range("J15").Value = 0.03119957432329093
$ws.Range("O15").Value = 0.2396231127748354
$ws.Range("P15").Value = 0.2396231127748355
$ws.Range("Q15").Value = 1.573212621782222
$ws.Range("R15").Value = 14.15891359604
$ws.Range("S15").Value = 0.007476139116596801
$ws.Range("T15").Value = 0.007476139116596803
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3599813333333333
$ws.Range("H16").Value = 1.079944
$ws.Range("I16").Value = 0.03119957432329092
$ws.Range("J16").Value = 0.03119957432329093
$ws.Range("M16").Value = 6.452372666666666
$ws.Range("N16").Value = 19.357118
$ws.Range("O16").Value = 0.3537860524377295
$ws.Range("P16").Value = 0.3537860524377295
$ws.Range("Q16").Value = 2.322733715710222
$ws.Range("R16").Value = 20.904603441392
$ws.Range("S16").Value = 0.01103797423757464
$ws.Range("T16").Value = 0.01103797423757464
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3599813333333333
$ws.Range("H17").Value = 1.079944
$ws.Range("I17").Value = 0.03119957432329092
$ws.Range("J17").Value = 0.03119957432329093
$ws.Range("M17").Value = 1.570781
$ws.Range("N17").Value = 4.712343
$ws.Range("O17").Value = 0.08612652088511148
$ws.Range("P17").Value = 0.0861265208851115
$ws.Range("Q17").Value = 0.5654518387546666
$ws.Range("R17").Value = 5.089066548792
$ws.Range("S17").Value = 0.002687110789561504
$ws.Range("T17").Value = 0.002687110789561504
